## chore: publish terminology IG 2.0.0 (#33)
## 1. Update Metadata sheet (Version + Date)
## 2. Rename "Concepts" -> "Properties"; repurpose it to hold FHIR property
##    definitions (status / effectiveDate)
## 3. Re-create a fresh "Concepts" sheet (copied from the original, so the
##    existing header/data rows + styles are preserved) positioned after
##    "Properties"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet: bump Version and Date
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "1.0.1"

# "2025-09-22" looks like a date to the type-inferencer, which would turn
# it into a date-serial cell; route it through a text formula + paste-values
# so it lands back as a literal string (matching the original cell's type
# and style).
$meta.Range("B8").Formula = "=""2025-09-22"""
$meta.Range("B8").Copy()
$meta.Range("B8").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2) Duplicate the existing "Concepts" sheet so the copy keeps all of the
#    original concept rows + formatting; this duplicate will become the
#    new "Concepts" tab.
# ---------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Copy($null, $concepts)
$newConcepts = $wb.Worksheets.Item("Concepts (2)")

# ---------------------------------------------------------------------
# 3) Turn the original sheet into "Properties": drop the old concept rows
#    (4-20) and overwrite the first three rows with the property table.
# ---------------------------------------------------------------------
$concepts.Range("A4:D20").EntireRow.Delete()

$concepts.Range("A1").Value = "Code"
$concepts.Range("B1").Value = "Uri"
$concepts.Range("C1").Value = "Description"
$concepts.Range("D1").Value = "Type"

$concepts.Range("A2").Value = "status"
$concepts.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$concepts.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$concepts.Range("D2").Value = "code"

$concepts.Range("A3").Value = "effectiveDate"
$concepts.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$concepts.Range("C3").Value = "The date at which the concept status was last changed."
$concepts.Range("D3").Value = "dateTime"

$concepts.Name = "Properties"

# ---------------------------------------------------------------------
# 4) Rename the duplicated sheet back to "Concepts" (it already sits right
#    after "Properties" because of the Copy(After:=) call above).
# ---------------------------------------------------------------------
$newConcepts.Name = "Concepts"

# Keep the original active tab (Metadata was the selected sheet before
# the edit, and the sheet-restructuring above leaves "Concepts" selected).
$meta.Activate()
